$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.788.25'
$ws.Range("E2").Value = '  +0.69%  '

$ws.Range("E3").Value = '  +0.55%  '

$ws.Range("E4").Value = '  +0.52%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '216.60'
$ws.Range("E5").Value = '  +1.24%  '

$ws.Range("E6").Value = '  +0.60%  '

$ws.Range("E7").Value = '  +0.54%  '

$ws.Range("E8").Value = '  +0.49%  '

$ws.Range("E9").Value = '  +0.00%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.25'
$ws.Range("E10").Value = '  +1.18%  '

$ws.Range("E11").Value = '  +0.20%  '

$ws.Range("D12").Value = '1.872.98'
$ws.Range("E12").Value = '  +0.37%  '

$ws.Range("D13").Value = '1.667.24'
$ws.Range("E13").Value = '  +1.51%  '

$ws.Range("E14").Value = '  +0.82%  '

$ws.Range("E15").Value = '  +0.85%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.49'
$ws.Range("E16").Value = '  -0.03%  '

$ws.Range("D17").Value = '26.792.42'
$ws.Range("E17").Value = '  +0.67%  '

$ws.Range("E18").Value = '  +0.03%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '217.85'
$ws.Range("E19").Value = '  +0.18%  '

$ws.Range("E20").Value = '  +0.47%  '

$ws.Range("E21").Value = '  +15.90%  '

$ws.Range("E22").Value = '  +0.85%  '

$ws.Range("E23").Value = '  +0.47%  '

$ws.Range("E24").Value = '  +0.60%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.66'
$ws.Range("E25").Value = '  -1.28%  '

$ws.Range("E26").Value = '  +0.52%  '

$ws.Range("E27").Value = '  +0.09%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.21'
$ws.Range("E28").Value = '  +4.06%  '

$ws.Range("E29").Value = '  +0.75%  '

$ws.Range("E30").Value = '  +0.45%  '

$ws.Range("E31").Value = '  +1.18%  '

$ws.Range("E32").Value = '  -0.65%  '

$ws.Range("E33").Value = '  +0.41%  '

$ws.Range("B34").Value = 'LidoDAOToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.55'
$ws.Range("E34").Value = '  +2.17%  '

$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '1.279.84'
$ws.Range("E35").Value = '  +0.63%  '

$ws.Range("E36").Value = '  +2.19%  '

$ws.Range("E37").Value = '  +1.65%  '

$ws.Range("E38").Value = '  +5.51%  '

$ws.Range("E39").Value = '  +3.54%  '

$ws.Range("E40").Value = '  +0.48%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.818'
$ws.Range("E41").Value = '  +2.01%  '

$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.45'

$ws.Range("D44").Value = '1.799.97'
$ws.Range("E44").Value = '  +1.43%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.03'
$ws.Range("E45").Value = '  -1.80%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '59.74'
$ws.Range("E46").Value = '  +7.96%  '

$ws.Range("E47").Value = '  +0.85%  '

$ws.Range("E48").Value = '  +0.62%  '

$ws.Range("E49").Value = '  +0.48%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.82'
$ws.Range("E50").Value = '  +2.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0980'
$ws.Range("E51").Value = '  +1.61%  '
